# Update the "K" column (G) values for the save_data sheet.
# The K column previously held "Strike#"-derived counts; this regenerates
# the column using the new K definition (std/mean recalculated, s_vals
# calculated and written). Only the G column values change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 3
    5  = 2
    6  = 4
    7  = 1
    8  = 2
    9  = 0
    10 = 2
    11 = 2
    12 = 1
    13 = 0
    14 = 1
    15 = 1
    16 = 1
    17 = 0
    18 = 0
    19 = 0
    20 = 1
    21 = 2
    22 = 1
    23 = 1
    24 = 0
    25 = 0
    26 = 1
    27 = 0
    28 = 0
    29 = 1
    30 = 2
    31 = 1
    32 = 1
    33 = 1
    34 = 1
    35 = 3
    36 = 1
    37 = 1
    38 = 1
    39 = 2
    40 = 0
    41 = 0
    42 = 0
    43 = 1
    44 = 0
    45 = 1
    46 = 1
    47 = 1
    48 = 1
    49 = 0
    50 = 2
    51 = 3
    52 = 0
    54 = 3
    55 = 0
    56 = 1
    57 = 2
    58 = 0
    59 = 1
    60 = 2
    61 = 1
    62 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
